$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 190, shifting rows 190:249 down to 191:250
$ws.Rows.Item(190).Insert()

# Populate the new row 190 - same descriptive columns as the old row190 (now row191), new date/price stats
$ws.Range("A190").Value = 4
$ws.Range("B190").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C190").Value = "Los Lagos"
$ws.Range("D190").Value = 44642
$ws.Range("D190").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E190").Value = 10
$ws.Range("F190").Value = "Fruta"
$ws.Range("G190").Value = 100102
$ws.Range("H190").Value = "Cítricos"
$ws.Range("I190").Value = 100102006
$ws.Range("J190").Value = "Pomelo"
$ws.Range("K190").Value = "Start Ruby"
$ws.Range("L190").Value = "Primera"
$ws.Range("M190").Value = 200
$ws.Range("N190").Value = 12000
$ws.Range("O190").Value = 13000
$ws.Range("P190").Value = 12500
$ws.Range("Q190").Value = '$/caja 14 kilos empedrada'
$ws.Range("R190").Value = "Región de O'Higgins"
$ws.Range("S190").Value = 893
$ws.Range("T190").Value = 14
